# Generate Report for Handoff
# Splits the existing "ab455dcf" entry's slot into two new entries
# (24eb05da, 5959d949) and pushes the original ab455dcf entry down to
# make room, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, $rowNum, $values) {
    foreach ($col in $values.Keys) {
        $colIndex = [int][char]$col - [int][char]'A' + 1
        $ws.Cells.Item($rowNum, $colIndex).Value = $values[$col]
    }
}

# ---------------------------------------------------------------------
# Overview sheet (7 columns, hyperlinks in column B)
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Rows.Item(4).Insert()
$wsOv.Rows.Item(4).Insert()

Set-RowValues $wsOv 3 @{
    A = "24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    B = "e2e\24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    C = ".md"
    D = ""
    E = "Ready for handoff"
    F = "Ready for handoff"
    G = "2016-08-17 18:39:49"
}
Set-RowValues $wsOv 4 @{
    A = "5959d949-8865-40ce-9726-2de26f357b70.md"
    B = "e2e\5959d949-8865-40ce-9726-2de26f357b70.md"
    C = ".md"
    D = ""
    E = "Ready for handoff"
    F = "Ready for handoff"
    G = "2016-08-17 18:39:49"
}
Set-RowValues $wsOv 5 @{
    A = "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md"
    B = "e2e\ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md"
    C = ".md"
    D = ""
    E = "Ready for handoff"
    F = "Ready for handoff"
    G = "2016-08-17 18:38:42"
}

foreach ($h in $wsOv.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$3') {
        $h.TextToDisplay = "e2e\24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    }
}
$null = $wsOv.Hyperlinks.Add($wsOv.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5652fc15c7bdf3d91a1f3972e2a412988431028c/e2e/5959d949-8865-40ce-9726-2de26f357b70.md", "", "", "e2e\5959d949-8865-40ce-9726-2de26f357b70.md")
$null = $wsOv.Hyperlinks.Add($wsOv.Cells.Item(5, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5652fc15c7bdf3d91a1f3972e2a412988431028c/e2e/ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md", "", "", "e2e\ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md")

foreach ($lo in $wsOv.ListObjects) {
    $lo.Resize($wsOv.Range("A1:G5"))
}

# ---------------------------------------------------------------------
# zh-cn sheet (16 columns, hyperlinks in columns A and I)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(4).Insert()
$wsZh.Rows.Item(4).Insert()

Set-RowValues $wsZh 3 @{
    A = "24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    B = ".md"
    C = "Ready for handoff"
    D = "e2e"
    E = "ht"
    F = "False"
    G = "24eb05da-1341-4b9a-a5ce-063982df19d5.2d6cb47790e59caeb9a2670e1c8c820bbc991e03.zh-cn.xlf"
    H = "2016-08-17 18:39:44"
    I = ""
    J = ""
    K = "0001-01-01 00:00:00"
    L = ""
    M = "True"
    N = ""
    O = "False"
    P = ""
}
Set-RowValues $wsZh 4 @{
    A = "5959d949-8865-40ce-9726-2de26f357b70.md"
    B = ".md"
    C = "Ready for handoff"
    D = "e2e"
    E = "ht"
    F = "False"
    G = "5959d949-8865-40ce-9726-2de26f357b70.f8d65483db9e32e8074e405776ba9f79fc23652e.zh-cn.xlf"
    H = "2016-08-17 18:39:44"
    I = ""
    J = ""
    K = "0001-01-01 00:00:00"
    L = ""
    M = "True"
    N = ""
    O = "False"
    P = ""
}
Set-RowValues $wsZh 5 @{
    A = "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md"
    B = ".md"
    C = "Ready for handoff"
    D = "e2e"
    E = "ht"
    F = "False"
    G = "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.9dbf2f61a274599fffb6b88e6dd99e92bb2ba384.zh-cn.xlf"
    H = "2016-08-17 18:38:37"
    I = ""
    J = ""
    K = "0001-01-01 00:00:00"
    L = ""
    M = "True"
    N = ""
    O = "False"
    P = ""
}

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$3') {
        $h.TextToDisplay = "24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    }
}
$null = $wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5652fc15c7bdf3d91a1f3972e2a412988431028c/e2e/5959d949-8865-40ce-9726-2de26f357b70.md", "", "", "5959d949-8865-40ce-9726-2de26f357b70.md")
$null = $wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5652fc15c7bdf3d91a1f3972e2a412988431028c/e2e/ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md", "", "", "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md")

foreach ($lo in $wsZh.ListObjects) {
    $lo.Resize($wsZh.Range("A1:P5"))
}

# ---------------------------------------------------------------------
# de-de sheet (16 columns, hyperlinks in columns A and I)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(4).Insert()
$wsDe.Rows.Item(4).Insert()

Set-RowValues $wsDe 3 @{
    A = "24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    B = ".md"
    C = "Ready for handoff"
    D = "e2e"
    E = "ht"
    F = "False"
    G = "24eb05da-1341-4b9a-a5ce-063982df19d5.2d6cb47790e59caeb9a2670e1c8c820bbc991e03.de-de.xlf"
    H = "2016-08-17 18:39:49"
    I = ""
    J = ""
    K = "0001-01-01 00:00:00"
    L = ""
    M = "True"
    N = ""
    O = "False"
    P = ""
}
Set-RowValues $wsDe 4 @{
    A = "5959d949-8865-40ce-9726-2de26f357b70.md"
    B = ".md"
    C = "Ready for handoff"
    D = "e2e"
    E = "ht"
    F = "False"
    G = "5959d949-8865-40ce-9726-2de26f357b70.f8d65483db9e32e8074e405776ba9f79fc23652e.de-de.xlf"
    H = "2016-08-17 18:39:49"
    I = ""
    J = ""
    K = "0001-01-01 00:00:00"
    L = ""
    M = "True"
    N = ""
    O = "False"
    P = ""
}
Set-RowValues $wsDe 5 @{
    A = "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md"
    B = ".md"
    C = "Ready for handoff"
    D = "e2e"
    E = "ht"
    F = "False"
    G = "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.9dbf2f61a274599fffb6b88e6dd99e92bb2ba384.de-de.xlf"
    H = "2016-08-17 18:38:42"
    I = ""
    J = ""
    K = "0001-01-01 00:00:00"
    L = ""
    M = "True"
    N = ""
    O = "False"
    P = ""
}

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$3') {
        $h.TextToDisplay = "24eb05da-1341-4b9a-a5ce-063982df19d5.md"
    }
}
$null = $wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5652fc15c7bdf3d91a1f3972e2a412988431028c/e2e/5959d949-8865-40ce-9726-2de26f357b70.md", "", "", "5959d949-8865-40ce-9726-2de26f357b70.md")
$null = $wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5652fc15c7bdf3d91a1f3972e2a412988431028c/e2e/ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md", "", "", "ab455dcf-2a9e-4d52-a91e-51eb4080f05d.md")

foreach ($lo in $wsDe.ListObjects) {
    $lo.Resize($wsDe.Range("A1:P5"))
}

Write-Output "done"
